$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2222.6
$ws.Range("J17").Value = 2153.25
$ws.Range("L17").Value = 6459.75
$ws.Range("N17").Value = -6795.75

$ws.Range("H98").Value = 1073.6086
$ws.Range("I98").Value = 854.2273
$ws.Range("K98").Value = 854.2273
$ws.Range("M98").Value = 643.7727

$ws.Range("H111").Value = 909.2
$ws.Range("I111").Value = 909.2
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2727.6
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 339.3999999999996
$ws.Range("N111").ClearContents()

$ws.Range("H116").Value = 5983
$ws.Range("I116").Value = 5846
$ws.Range("J116").Value = 6165.6665
$ws.Range("K116").Value = 5846
$ws.Range("L116").Value = 6165.6665
$ws.Range("M116").Value = -2404
$ws.Range("N116").Value = -13049.6665

$ws.Range("H122").Value = 1073.6086
$ws.Range("I122").Value = 854.2273
$ws.Range("K122").Value = 2562.6819
$ws.Range("M122").Value = -112.6819

$ws.Range("H137").Value = 5023
$ws.Range("I137").Value = 899.75
$ws.Range("K137").Value = 2699.25
$ws.Range("M137").Value = -149.25

$ws.Range("H141").Value = 5164.2354
$ws.Range("I141").Value = 5559.467
$ws.Range("K141").Value = 16678.401
$ws.Range("M141").Value = -11498.401

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 54985.547
$ws.Range("I32").Value = 66549.234
$ws.Range("J32").Value = 15669
$ws.Range("K32").Value = 66549.234
$ws.Range("L32").Value = 15669
$ws.Range("M32").Value = -66262.234
$ws.Range("N32").Value = -16243

$ws.Range("H88").Value = 1587.9474
$ws.Range("I88").Value = 1149.6
$ws.Range("J88").Value = 1744.5
$ws.Range("K88").Value = 1149.6
$ws.Range("L88").Value = 1744.5
$ws.Range("M88").Value = -743.5999999999999
$ws.Range("N88").Value = -2556.5

$ws.Range("H91").Value = 1587.9474
$ws.Range("I91").Value = 1149.6
$ws.Range("J91").Value = 1744.5
$ws.Range("K91").Value = 1149.6
$ws.Range("L91").Value = 1744.5
$ws.Range("M91").Value = 254.4000000000001
$ws.Range("N91").Value = -4552.5

$ws.Range("H122").Value = 2389.889
$ws.Range("I122").Value = 2168.5
$ws.Range("J122").Value = 2832.6667
$ws.Range("K122").Value = 6505.5
$ws.Range("L122").Value = 8498.000100000001
$ws.Range("M122").Value = -4055.5
$ws.Range("N122").Value = -13398.0001

$ws.Range("H132").Value = 205137.2
$ws.Range("I132").Value = 205137.2
$ws.Range("K132").Value = 615411.6000000001
$ws.Range("M132").Value = -612881.6000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2345.7058
$ws.Range("I86").Value = 2255.75
$ws.Range("K86").Value = 2255.75
$ws.Range("M86").Value = -1132.75

$ws.Range("H89").Value = 2345.7058
$ws.Range("I89").Value = 2255.75
$ws.Range("K89").Value = 11278.75
$ws.Range("M89").Value = -5662.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 49138.59
$ws.Range("I134").Value = 55639.42
$ws.Range("K134").Value = 166918.26
$ws.Range("M134").Value = -164383.26

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 152.66667
$ws.Range("I7").Value = 72.666664
$ws.Range("K7").Value = 217.999992
$ws.Range("M7").Value = -105.999992

$ws.Range("H131").Value = 1790280.8
$ws.Range("I131").Value = 1867.5385
$ws.Range("J131").Value = 2330964
$ws.Range("K131").Value = 5602.6155
$ws.Range("L131").Value = 6992892
$ws.Range("M131").Value = -562.6154999999999
$ws.Range("N131").Value = -7002972

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7966.1665
$ws.Range("I70").Value = 11504.5
$ws.Range("J70").Value = 6197
$ws.Range("K70").Value = 11504.5
$ws.Range("L70").Value = 6197
$ws.Range("M70").Value = -11234.5
$ws.Range("N70").Value = -6737

$ws.Range("H73").Value = 7966.1665
$ws.Range("I73").Value = 11504.5
$ws.Range("J73").Value = 6197
$ws.Range("K73").Value = 11504.5
$ws.Range("L73").Value = 6197
$ws.Range("M73").Value = -10568.5
$ws.Range("N73").Value = -8069

$ws.Range("H102").Value = 1983.7407
$ws.Range("I102").Value = 1415.7391
$ws.Range("K102").Value = 1415.7391
$ws.Range("M102").Value = 206.2609

$ws.Range("H107").Value = 48603.953
$ws.Range("I107").Value = 63450
$ws.Range("J107").Value = 1096.6
$ws.Range("K107").Value = 63450
$ws.Range("L107").Value = 1096.6
$ws.Range("M107").Value = -61530
$ws.Range("N107").Value = -4936.6

$ws.Range("H113").Value = 3399
$ws.Range("I113").Value = 3298
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 3298
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -1128
$ws.Range("N113").Value = -7840

$ws.Range("H122").Value = 1997.6111
$ws.Range("I122").Value = 1969.4482
$ws.Range("J122").Value = 2114.2856
$ws.Range("K122").Value = 5908.3446
$ws.Range("L122").Value = 6342.8568
$ws.Range("M122").Value = -3458.3446
$ws.Range("N122").Value = -11242.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3898.0527
$ws.Range("I7").Value = 3851.2307
$ws.Range("K7").Value = 3851.2307
$ws.Range("M7").Value = -3739.2307

$ws.Range("H40").Value = 3886.5
$ws.Range("I40").Value = 3188.4
$ws.Range("J40").Value = 5050
$ws.Range("K40").Value = 3188.4
$ws.Range("L40").Value = 5050
$ws.Range("M40").Value = -3052.4
$ws.Range("N40").Value = -5322

$ws.Range("H93").Value = 1577.05
$ws.Range("I93").Value = 945.9286
$ws.Range("K93").Value = 945.9286
$ws.Range("M93").Value = 302.0714

$ws.Range("H100").Value = 2409.0476
$ws.Range("I100").Value = 2142.7144
$ws.Range("J100").Value = 2941.7144
$ws.Range("K100").Value = 2142.7144
$ws.Range("L100").Value = 2941.7144
$ws.Range("M100").Value = -1601.7144
$ws.Range("N100").Value = -4023.7144

$ws.Range("H122").Value = 3837.3572
$ws.Range("I122").Value = 3354.875
$ws.Range("J122").Value = 4030.35
$ws.Range("K122").Value = 10064.625
$ws.Range("L122").Value = 12091.05
$ws.Range("M122").Value = -7614.625
$ws.Range("N122").Value = -16991.05

$ws.Range("H126").Value = 3898.0527
$ws.Range("I126").Value = 3851.2307
$ws.Range("K126").Value = 11553.6921
$ws.Range("M126").Value = -9083.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 569.5
$ws.Range("I107").Value = 639.625
$ws.Range("J107").Value = 289
$ws.Range("K107").Value = 1918.875
$ws.Range("L107").Value = 867
$ws.Range("M107").Value = 1.125
$ws.Range("N107").Value = -4707

$ws.Range("H122").Value = 1266.4615
$ws.Range("I122").Value = 1288.25
$ws.Range("J122").Value = 1005
$ws.Range("K122").Value = 3864.75
$ws.Range("L122").Value = 3015
$ws.Range("M122").Value = -1414.75
$ws.Range("N122").Value = -7915

$ws.Range("H136").Value = 2179.7693
$ws.Range("I136").Value = 1862.3871
$ws.Range("K136").Value = 5587.1613
$ws.Range("M136").Value = -3037.1613
